$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Passengers row: add "n/a" in B11
$ws.Range("B11").Value = "n/a"

# System designer row: update B12 message
$ws.Range("B12").Value = "no message, but creation of functionality -  variable and default capacity"

# Row 12 grows taller to fit the wrapped text (matches autofit behavior in Excel)
$ws.Rows.Item(12).RowHeight = 32

# Update the active selection as recorded in the saved file
$ws.Range("A15").Select()
